$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into a cell that lives on a protected worksheet
# without disturbing the sheet's protection settings (password etc.) and
# without leaving the cell's visible number format changed.
# ---------------------------------------------------------------------------
function Set-ProtectedValue($addr, $val) {
    $r = $ws.Range($addr)
    $fmt = $r.NumberFormat
    $r.Locked = $false
    $r.Value = $val
    $r.Style = "Normal"
    $r.NumberFormat = $fmt
}

# ---------------------------------------------------------------------------
# Update the "as of" date in the confidentiality footnote (row 18, col A):
#   2021-06-09  ->  2021-06-10
# Edit only the two digits that changed, in place, rather than rewriting
# the whole string.
# ---------------------------------------------------------------------------
$footnote = $ws.Range("A18")
$footnote.Locked = $false
$chars = $footnote.Characters(122, 2)
$chars.Text = "10"
$footnote.Style = "Normal"

# ---------------------------------------------------------------------------
# Updated Weight (D) / Percent Change (E) figures, rows 2-15
# ---------------------------------------------------------------------------
Set-ProtectedValue "D2"  0.05742459586524589
Set-ProtectedValue "E2"  0.004421148587055557

Set-ProtectedValue "D3"  0.02078827951958594
Set-ProtectedValue "E3"  0.003210576015108435

Set-ProtectedValue "D4"  0.02853955302254518
Set-ProtectedValue "E4"  0.006969919295671145

Set-ProtectedValue "D5"  0.0307459766467203
Set-ProtectedValue "E5"  -0.01555716353111436

Set-ProtectedValue "D6"  0.03109408381399823
Set-ProtectedValue "E6"  -0.01999158249158262

Set-ProtectedValue "D7"  0.01906344776593065
Set-ProtectedValue "E7"  -0.01681883709754917

Set-ProtectedValue "D8"  0.01024625964737791
Set-ProtectedValue "E8"  -0.001341081805990241

Set-ProtectedValue "D9"  0.01037843567705863
Set-ProtectedValue "E9"  -0.000945715906941369

Set-ProtectedValue "D10" 0.07000683781935725
Set-ProtectedValue "E10" 0.004794885455514208

Set-ProtectedValue "D11" 0.07011872940883945
Set-ProtectedValue "E11" 0.00478723404255299

Set-ProtectedValue "D12" 0.1483407655085996
Set-ProtectedValue "E12" 0.005999011927447118

Set-ProtectedValue "D13" 0.3925523387632381
Set-ProtectedValue "E13" 0.003553783479240735

Set-ProtectedValue "D14" 0.110700696541503
Set-ProtectedValue "E14" 0.01137250265988898

Set-ProtectedValue "E15" 0.003290659669625606
